$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; unprotect to allow writing into locked cells,
# then restore protection once all edits are applied.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A59).
$newDisclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."
$ws.Range("A59").Value = $newDisclaimer

# Update Weight (D) and Percent Change (E) values for rows 2-56.
$ws.Range("D2").Value = 0.01981813486156674
$ws.Range("E2").Value = -0.01603125420988816
$ws.Range("D3").Value = 0.01931887698481705
$ws.Range("E3").Value = -0.009765984890363177
$ws.Range("D4").Value = 0.01891733482084296
$ws.Range("E4").Value = -0.03183922321327759
$ws.Range("D5").Value = 0.02014545580108285
$ws.Range("E5").Value = -0.007103477523324631
$ws.Range("D6").Value = 0.01925480110758714
$ws.Range("E6").Value = -0.01830282861896837
$ws.Range("D7").Value = 0.008000585226345367
$ws.Range("E7").Value = -0.0002224694104561165
$ws.Range("D8").Value = 0.01984643504067661
$ws.Range("E8").Value = -0.02000825082508251
$ws.Range("D9").Value = 0.02396437808398532
$ws.Range("E9").Value = -0.002005347593582951
$ws.Range("D10").Value = 0.02063296643367373
$ws.Range("E10").Value = -0.01891772987241536
$ws.Range("D11").Value = 0.01935038095778842
$ws.Range("E11").Value = 0.01757774773034559
$ws.Range("D12").Value = 0.01365528139188469
$ws.Range("E12").Value = 0.02150677789363931
$ws.Range("D13").Value = 0.01379553636759904
$ws.Range("E13").Value = 0.002322340919646892
$ws.Range("D14").Value = 0.008582607777850358
$ws.Range("E14").Value = -0.00240564081294059
$ws.Range("D15").Value = 0.01444697445276976
$ws.Range("E15").Value = 0.001478415138970846
$ws.Range("D16").Value = 0.02127034342356344
$ws.Range("E16").Value = 0.0362665684830632
$ws.Range("D17").Value = 0.02367247686549353
$ws.Range("E17").Value = 0.009263157894736862
$ws.Range("D18").Value = 0.02162988029024237
$ws.Range("E18").Value = -0.006879299562226326
$ws.Range("D19").Value = 0.0193477111295705
$ws.Range("E19").Value = -0.01032179720704318
$ws.Range("D20").Value = 0.0195805201501725
$ws.Range("E20").Value = -0.01381692573402404
$ws.Range("D21").Value = 0.02551661620987982
$ws.Range("E21").Value = -0.002253053480374523
$ws.Range("D22").Value = 0.01684554812374262
$ws.Range("E22").Value = 0.01388360593381521
$ws.Range("D23").Value = 0.02111086568468012
$ws.Range("E23").Value = -0.01392823418319145
$ws.Range("D24").Value = 0.0190999510709482
$ws.Range("E24").Value = -0.008163265306122436
$ws.Range("D25").Value = 0.02020259012494618
$ws.Range("E25").Value = -0.008308004052684947
$ws.Range("D26").Value = 0.01843783367290582
$ws.Range("E26").Value = -0.01766579785693612
$ws.Range("D27").Value = 0.01991763045982096
$ws.Range("E27").Value = -0.01835501858736055
$ws.Range("D28").Value = 0.02131590849181583
$ws.Range("E28").Value = -0.01102204408817653
$ws.Range("D29").Value = 0.02011911349599944
$ws.Range("E29").Value = -0.01021798365122617
$ws.Range("D30").Value = 0.01997031506998777
$ws.Range("E30").Value = 0.008333333333333304
$ws.Range("D31").Value = 0.02164678920228915
$ws.Range("E31").Value = -0.02905795969379799
$ws.Range("D32").Value = 0.02251661923568514
$ws.Range("E32").Value = -0.005438477226376692
$ws.Range("D33").Value = 0.01995465207777601
$ws.Range("E33").Value = -0.01461038961038952
$ws.Range("D34").Value = 0.02007888808418289
$ws.Range("E34").Value = -0.02210796915167079
$ws.Range("D35").Value = 0.02013050476306254
$ws.Range("E35").Value = 0.000212201591512029
$ws.Range("D36").Value = 0.01710825922038524
$ws.Range("E36").Value = -0.01685393258426959
$ws.Range("D37").Value = 0.02090030523256073
$ws.Range("E37").Value = -0.01039812646370031
$ws.Range("D38").Value = 0.01969852655740425
$ws.Range("E38").Value = 0.01268602098072691
$ws.Range("D39").Value = 0.01981760089592316
$ws.Range("E39").Value = -0.01923802338740099
$ws.Range("D40").Value = 0.01688719744394206
$ws.Range("E40").Value = 0.00512236767216856
$ws.Range("D41").Value = 0.01307859849681552
$ws.Range("E41").Value = -0.009363091997822481
$ws.Range("D42").Value = 0.0168533796198485
$ws.Range("E42").Value = 0.01394052044609673
$ws.Range("D43").Value = 0.01924109598940185
$ws.Range("E43").Value = -0.02170152539707504
$ws.Range("D44").Value = 0.01291823081519289
$ws.Range("E44").Value = -0.0005373455131647553
$ws.Range("D45").Value = 0.01692813480995006
$ws.Range("E45").Value = -0.04889178617992174
$ws.Range("D46").Value = 0.01689467296295222
$ws.Range("E46").Value = 0.01327433628318597
$ws.Range("D47").Value = 0.01369265898693546
$ws.Range("E47").Value = -0.002547770700636831
$ws.Range("D48").Value = 0.0209777302508802
$ws.Range("E48").Value = -0.01637536059731859
$ws.Range("D49").Value = 0.01928345726379273
$ws.Range("E49").Value = -0.02086929232700441
$ws.Range("D50").Value = 0.01800710138708255
$ws.Range("E50").Value = -0.008609271523178696
$ws.Range("D51").Value = 0.01929004284006359
$ws.Range("E51").Value = -0.01013120744062446
$ws.Range("D52").Value = 0.006287979418828235
$ws.Range("E52").Value = -0.007472826086956652
$ws.Range("D53").Value = 0.02150956003188843
$ws.Range("E53").Value = -0.008928571428571397
$ws.Range("D54").Value = 0.01931496123676411
$ws.Range("E54").Value = -0.0323356493853556
$ws.Range("D55").Value = 0.01921599960415347
$ws.Range("E55").Value = -0.01981252662973998
$ws.Range("E56").Value = -0.008134254625788651

$ws.Protect()
